$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("E2").Value = [double]"23.86000000000029"
$ws.Range("H2").Value = [double]"1.277954560719605e-16"
$ws.Range("K2").Value = [double]"46.13676814161559"
$ws.Range("L2").Value = "[41.85131668419478, 50.4222195990364]"
$ws.Range("O2").Value = [double]"1.641552918091964"
$ws.Range("P2").Value = "[1.54092132158058, 1.7421845146033483]"
$ws.Range("S2").Value = [double]"58.06247306845448"
$ws.Range("T2").Value = "[55.22427532833242, 60.90067080857654]"
$ws.Range("W2").Value = [double]"17.62630630630652"
$ws.Range("X2").Value = [double]"17.24416416416438"
$ws.Range("Y2").Value = [double]"18.00844844844866"

# Row 3 updates
$ws.Range("E3").Value = [double]"23.24000000000019"
$ws.Range("H3").Value = [double]"1.277954560719605e-16"
$ws.Range("K3").Value = [double]"48.58817404087332"
$ws.Range("L3").Value = "[41.374335709612815, 55.802012372133824]"
$ws.Range("O3").Value = [double]"-0.3270526886620004"
$ws.Range("P3").Value = "[-0.49057903299300065, -0.16352634433100022]"
$ws.Range("Q3").Value = [double]"0.0001123143493098056"
$ws.Range("R3").Value = [double]"0.0001123143493098056"
$ws.Range("S3").Value = [double]"53.46007092299402"
$ws.Range("T3").Value = "[49.383141414729664, 57.53700043125838]"
$ws.Range("W3").Value = [double]"1.2096896896897"
$ws.Range("X3").Value = [double]"0.6048448448448502"
$ws.Range("Y3").Value = [double]"1.81453453453455"
